$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column C (Förändrad) starting from the bottom.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

# The "Förändrad" date for every data row was bumped by one day (46081 -> 46082).
$ws.Range("C2:C$lastRow").Value = 46082
